$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.140.55'
$ws.Range('E2').Value = '  +3.42%  '
$ws.Range('D3').Value = '3.962.93'
$ws.Range('E3').Value = '  +3.94%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '''482.22'
$ws.Range('E5').Value = '  +7.71%  '
$ws.Range('D6').Value = '''148.65'
$ws.Range('E6').Value = '  +1.70%  '
$ws.Range('D7').Value = '''0.624'
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '''0.729'
$ws.Range('E9').Value = '  -0.94%  '
$ws.Range('D10').Value = '''0.168'
$ws.Range('E10').Value = '  +8.56%  '
$ws.Range('D11').Value = '''0.0000351'
$ws.Range('E11').Value = '  +8.55%  '
$ws.Range('D12').Value = '''42.91'
$ws.Range('E12').Value = '  -1.04%  '
$ws.Range('D13').Value = '4.589.10'
$ws.Range('E13').Value = '  +3.80%  '
$ws.Range('D14').Value = '''10.35'
$ws.Range('E14').Value = '  +0.17%  '
$ws.Range('D15').Value = '''14.68'
$ws.Range('E15').Value = '  -2.91%  '
$ws.Range('D16').Value = '3.947.98'
$ws.Range('E16').Value = '  +3.81%  '
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').Value = '''19.77'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').Value = '''1.13'
$ws.Range('E19').Value = '  -0.80%  '
$ws.Range('D20').Value = '69.235.79'
$ws.Range('E20').Value = '  +3.41%  '
$ws.Range('D21').Value = '''438.59'
$ws.Range('E21').Value = '  +3.59%  '
$ws.Range('B22').Value = 'InternetComputer(DFINITY)'
$ws.Range('C22').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D22').Value = '''14.40'
$ws.Range('E22').Value = '  -1.74%  '
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').Value = '''3.30'
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('D24').Value = '''87.98'
$ws.Range('E24').Value = '  +1.69%  '
$ws.Range('D25').Value = '''3.70'
$ws.Range('E25').Value = '  +7.83%  '
$ws.Range('D26').Value = '''38.44'
$ws.Range('E26').Value = '  +3.12%  '
$ws.Range('B27').Value = 'Filecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D27').Value = '''9.97'
$ws.Range('E27').Value = '  +2.95%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').Value = '''9.67'
$ws.Range('E28').Value = '  +1.83%  '
$ws.Range('D29').Value = '''729.49'
$ws.Range('E29').Value = '  -2.47%  '
$ws.Range('D30').Value = '''13.18'
$ws.Range('E30').Value = '  -3.65%  '
$ws.Range('D31').Value = '''0.127'
$ws.Range('E31').Value = '  -4.66%  '
$ws.Range('D32').Value = '''2.82'
$ws.Range('E32').Value = '  +3.01%  '
$ws.Range('D33').Value = '''41.92'
$ws.Range('E33').Value = '  -2.44%  '
$ws.Range('D34').Value = '''59.78'
$ws.Range('E34').Value = '  +3.42%  '
$ws.Range('D35').Value = '0.0₃0846'
$ws.Range('E35').Value = '  +25.10%  '
$ws.Range('D36').Value = '''0.149'
$ws.Range('E36').Value = '  -3.53%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = '''5.35'
$ws.Range('E38').Value = '  -2.47%  '
$ws.Range('D39').Value = '''0.0472'
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('D40').Value = '''3.06'
$ws.Range('E40').Value = '  +7.09%  '
$ws.Range('E41').Value = '  +11.15%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').Value = '''0.141'
$ws.Range('E42').Value = '  +0.87%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = '''0.997'
$ws.Range('E43').Value = '  -0.51%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').Value = '''2.55'
$ws.Range('E44').Value = '  +2.44%  '
$ws.Range('D45').Value = '''0.332'
$ws.Range('E45').Value = '  -4.14%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '''3.26'
$ws.Range('E46').Value = '  +1.01%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '''2.16'
$ws.Range('E47').Value = '  +2.31%  '
$ws.Range('D48').Value = '''148.55'
$ws.Range('E48').Value = '  +1.71%  '
$ws.Range('D49').Value = '''3.37'
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('D50').Value = '''2.94'
$ws.Range('E50').Value = '  +3.43%  '
$ws.Range('D51').Value = '''24.97'
$ws.Range('E51').Value = '  -0.17%  '
